# Fruta / hortaliza, semanal
#
# A new weekly price observation is inserted at row 28 of the sheet.
# All existing data rows from 28 through 49 shift down by one (to 29-50),
# and the new row 28 receives its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 28, pushing rows 28-49 down to 29-50.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new observation.
$ws.Range("A28").Value = 5
$ws.Range("B28").Value = "Macroferia Regional de Talca"
$ws.Range("C28").Value = "Maule"
$ws.Range("D28").Value = 44554
$ws.Range("E28").Value = 7
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100101
$ws.Range("H28").Value = "Berries"
$ws.Range("I28").Value = 100101001
$ws.Range("J28").Value = "Arándano (blue)"
$ws.Range("K28").Value = "Sin especificar"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 100
$ws.Range("N28").Value = 3800
$ws.Range("O28").Value = 3800
$ws.Range("P28").Value = 3800
$ws.Range("Q28").Value = "$/bandeja 2 kilos"
$ws.Range("R28").Value = "Provincia de Linares"
$ws.Range("S28").Value = 1900
$ws.Range("T28").Value = 2
